$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4676999.5
$ws.Range("J17").Value = 5455871.5
$ws.Range("L17").Value = 16367614.5
$ws.Range("N17").Value = -16367950.5
$ws.Range("H137").Value = 3811.868
$ws.Range("I137").Value = 1485.027
$ws.Range("J137").Value = 9192.6875
$ws.Range("K137").Value = 4455.081
$ws.Range("L137").Value = 27578.0625
$ws.Range("M137").Value = -1905.081
$ws.Range("N137").Value = -32678.0625
$ws.Range("H138").Value = 3775.7742
$ws.Range("I138").Value = 3823.2666
$ws.Range("K138").Value = 11469.7998
$ws.Range("M138").Value = -6329.799800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2624.25
$ws.Range("I2").Value = 2472.7646
$ws.Range("K2").Value = 2472.7646
$ws.Range("M2").Value = -2359.7646
$ws.Range("H32").Value = 3997.6743
$ws.Range("I32").Value = 1938.4412
$ws.Range("J32").Value = 11777
$ws.Range("K32").Value = 1938.4412
$ws.Range("L32").Value = 11777
$ws.Range("M32").Value = -1651.4412
$ws.Range("N32").Value = -12351
$ws.Range("H61").Value = 4107.4653
$ws.Range("I61").Value = 3088.0908
$ws.Range("K61").Value = 3088.0908
$ws.Range("M61").Value = -2876.0908
$ws.Range("H109").Value = 28162.25
$ws.Range("J109").Value = 28162.25
$ws.Range("L109").Value = 28162.25
$ws.Range("N109").Value = -30936.25
$ws.Range("H116").Value = 2624.25
$ws.Range("I116").Value = 2472.7646
$ws.Range("K116").Value = 2472.7646
$ws.Range("M116").Value = -178.7646
$ws.Range("H122").Value = 1811.5264
$ws.Range("I122").Value = 1723.0625
$ws.Range("K122").Value = 5169.1875
$ws.Range("M122").Value = -2719.1875
$ws.Range("H132").Value = 2823.0334
$ws.Range("I132").Value = 2843.7585
$ws.Range("K132").Value = 8531.2755
$ws.Range("M132").Value = -6001.2755
$ws.Range("H136").Value = 4107.4653
$ws.Range("I136").Value = 3088.0908
$ws.Range("K136").Value = 9264.2724
$ws.Range("M136").Value = -6714.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2624.25
$ws.Range("I3").Value = 2472.7646
$ws.Range("K3").Value = 2472.7646
$ws.Range("M3").Value = -2358.7646
$ws.Range("H86").Value = 272290.56
$ws.Range("I86").Value = 372150.6
$ws.Range("K86").Value = 372150.6
$ws.Range("M86").Value = -371027.6
$ws.Range("H89").Value = 272290.56
$ws.Range("I89").Value = 372150.6
$ws.Range("K89").Value = 1860753
$ws.Range("M89").Value = -1855137
$ws.Range("H103").Value = 17606.834
$ws.Range("J103").Value = 17606.834
$ws.Range("L103").Value = 17606.834
$ws.Range("N103").Value = -19950.834
$ws.Range("H134").Value = 13869.5
$ws.Range("I134").Value = 3930.6667
$ws.Range("K134").Value = 11792.0001
$ws.Range("M134").Value = -9257.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2691.8394
$ws.Range("I31").Value = 1173.2778
$ws.Range("K31").Value = 1173.2778
$ws.Range("M31").Value = -878.2778000000001
$ws.Range("H34").Value = 2691.8394
$ws.Range("I34").Value = 1173.2778
$ws.Range("K34").Value = 1173.2778
$ws.Range("M34").Value = -971.2778000000001
$ws.Range("H132").Value = 4140.846
$ws.Range("I132").Value = 3652.5833
$ws.Range("K132").Value = 10957.7499
$ws.Range("M132").Value = -8427.749899999999
$ws.Range("H134").Value = 5926.963
$ws.Range("I134").Value = 6392.1816
$ws.Range("K134").Value = 19176.5448
$ws.Range("M134").Value = -16641.5448

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 163.26471
$ws.Range("I2").Value = 23.6
$ws.Range("J2").Value = 551.2222
$ws.Range("K2").Value = 141.6
$ws.Range("L2").Value = 3307.3332
$ws.Range("M2").Value = -28.60000000000002
$ws.Range("N2").Value = -3533.3332
$ws.Range("H121").Value = 655.8421
$ws.Range("I121").Value = 645.3333
$ws.Range("J121").Value = 695.25
$ws.Range("K121").Value = 1935.9999
$ws.Range("L121").Value = 2085.75
$ws.Range("M121").Value = -625.9999
$ws.Range("N121").Value = -4705.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8235.625
$ws.Range("I70").Value = 8365.5
$ws.Range("J70").Value = 8192.333000000001
$ws.Range("K70").Value = 8365.5
$ws.Range("L70").Value = 8192.333000000001
$ws.Range("M70").Value = -8095.5
$ws.Range("N70").Value = -8732.333000000001
$ws.Range("H73").Value = 8235.625
$ws.Range("I73").Value = 8365.5
$ws.Range("J73").Value = 8192.333000000001
$ws.Range("K73").Value = 8365.5
$ws.Range("L73").Value = 8192.333000000001
$ws.Range("M73").Value = -7429.5
$ws.Range("N73").Value = -10064.333
$ws.Range("H80").Value = 20565.8
$ws.Range("I80").Value = 13607.667
$ws.Range("J80").Value = 31003
$ws.Range("K80").Value = 13607.667
$ws.Range("L80").Value = 31003
$ws.Range("M80").Value = -12609.667
$ws.Range("N80").Value = -32999
$ws.Range("H83").Value = 20565.8
$ws.Range("I83").Value = 13607.667
$ws.Range("J83").Value = 31003
$ws.Range("K83").Value = 68038.33499999999
$ws.Range("L83").Value = 155015
$ws.Range("M83").Value = -63046.33499999999
$ws.Range("N83").Value = -164999
$ws.Range("H102").Value = 2157.9048
$ws.Range("I102").Value = 2180.4119
$ws.Range("K102").Value = 2180.4119
$ws.Range("M102").Value = -558.4119000000001
$ws.Range("H132").Value = 8785.612999999999
$ws.Range("I132").Value = 10433.72
$ws.Range("K132").Value = 31301.16
$ws.Range("M132").Value = -28771.16

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2364
$ws.Range("N40").ClearContents()
$ws.Range("H82").Value = 2865.5833
$ws.Range("I82").Value = 2324.3333
$ws.Range("J82").Value = 3406.8333
$ws.Range("K82").Value = 2324.3333
$ws.Range("L82").Value = 3406.8333
$ws.Range("M82").Value = -1963.3333
$ws.Range("N82").Value = -4128.8333
$ws.Range("H85").Value = 2865.5833
$ws.Range("I85").Value = 2324.3333
$ws.Range("J85").Value = 3406.8333
$ws.Range("K85").Value = 2324.3333
$ws.Range("L85").Value = 3406.8333
$ws.Range("M85").Value = -1076.3333
$ws.Range("N85").Value = -5902.8333
$ws.Range("H132").Value = 3446.2368
$ws.Range("I132").Value = 3064.963
$ws.Range("J132").Value = 4382.091
$ws.Range("K132").Value = 9194.889000000001
$ws.Range("L132").Value = 13146.273
$ws.Range("M132").Value = -6664.889000000001
$ws.Range("N132").Value = -18206.273
$ws.Range("H136").Value = 2517.7693
$ws.Range("I136").Value = 2217.1667
$ws.Range("K136").Value = 6651.500100000001
$ws.Range("M136").Value = -4101.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10026
$ws.Range("I32").Value = 10026
$ws.Range("K32").Value = 10026
$ws.Range("M32").Value = -9709
$ws.Range("H34").Value = 5500
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -797
$ws.Range("N34").Value = -10406
$ws.Range("H107").Value = 1330.0264
$ws.Range("I107").Value = 1211.0952
$ws.Range("K107").Value = 3633.2856
$ws.Range("M107").Value = -1713.2856
$ws.Range("H112").Value = 60387
$ws.Range("J112").Value = 60387
$ws.Range("L112").Value = 60387
$ws.Range("N112").Value = -63341
$ws.Range("H122").Value = 6161.625
$ws.Range("I122").Value = 3399.75
$ws.Range("J122").Value = 8923.5
$ws.Range("K122").Value = 10199.25
$ws.Range("L122").Value = 26770.5
$ws.Range("M122").Value = -7749.25
$ws.Range("N122").Value = -31670.5
$ws.Range("H132").Value = 15092.135
$ws.Range("I132").Value = 13950.538
$ws.Range("K132").Value = 41851.614
$ws.Range("M132").Value = -39321.614
$ws.Range("H136").Value = 1002.2692
$ws.Range("I136").Value = 1033.0416
$ws.Range("K136").Value = 3099.1248
$ws.Range("M136").Value = -549.1248000000001

